$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(2314,1012,1004,509,387,285,284,179,75,70,60,56,50,35,29,25,24,23,17,16,13,13,12,10,10,9,8,7,7,7,7,6,5,5,5,5,5,4,4,3,3,3,3,3,2,2,2,2,2,2,1,1,1,1,1,1,1,1,1,1,1,1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
